# Generate Report for handback
#
# The 9d491b91 file has now been handed back (in sync with en-US), so the
# localization-status report is regenerated:
#   - 9d491b91's row moves ahead of ae38855b's row on every sheet
#   - 9d491b91's status flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" and gains handback file/datetime data
#   - cb5280d3 (which depends on the 9d491b91 handoff) follows suit and is
#     now also "Handed back: in sync with en-US" with the full handback info

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"

$ws.Range("A3").Value = "ae38855b-7f4f-4934-8340-95b809897df5.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"

$ws.Range("A4").Value = "cb5280d3-793c-4a18-8880-5345aa846a9e.md"
$ws.Range("B4").Value = "Handed back: in sync with en-US"
$ws.Range("C4").Value = "Handed back: in sync with en-US"

$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("C5").Value = "Not to be localized"

# Rebuild the hyperlinks for rows 2-5 in the new order / targets.
$ws.Cells.Item(1, 1).Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5dc5f659cfefafe3787e5cc4efd5a0cb65c8d210/e2e/9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md", "", "", "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/481cee7df47d66e29231fd5af014ae7c5826663a/e2e/ae38855b-7f4f-4934-8340-95b809897df5.md", "", "", "ae38855b-7f4f-4934-8340-95b809897df5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5dc5f659cfefafe3787e5cc4efd5a0cb65c8d210/e2e/cb5280d3-793c-4a18-8880-5345aa846a9e.md", "", "", "cb5280d3-793c-4a18-8880-5345aa846a9e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/5dc5f659cfefafe3787e5cc4efd5a0cb65c8d210/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 2: 9d491b91 - now fully handed back
$ws.Range("A2").Value = "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.zh-cn.xlf"
$ws.Range("D2").Value = "2016-01-20 07:51:58"
$ws.Range("E2").Value = "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md"
$ws.Range("F2").Value = "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.zh-cn.xlf"
$ws.Range("G2").Value = "2016-01-20 07:52:43"
$ws.Range("H2").Value = "Include"

# Row 3: ae38855b - unchanged content, just shifted down a row
$ws.Range("A3").Value = "ae38855b-7f4f-4934-8340-95b809897df5.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.zh-cn.xlf"
$ws.Range("D3").Value = "2016-01-20 07:50:07"
$ws.Range("E3").Value = "ae38855b-7f4f-4934-8340-95b809897df5.md"
$ws.Range("F3").Value = "ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.zh-cn.xlf"
$ws.Range("G3").Value = "2016-01-20 07:50:54"
$ws.Range("H3").Value = "Include"

# Row 4: cb5280d3 dependency - now also handed back, mirroring 9d491b91's data
$ws.Range("A4").Value = "cb5280d3-793c-4a18-8880-5345aa846a9e.md"
$ws.Range("B4").Value = "Handed back: in sync with en-US"
$ws.Range("C4").Value = "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.zh-cn.xlf"
$ws.Range("D4").Value = "2016-01-20 07:51:58"
$ws.Range("E4").Value = "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md"
$ws.Range("F4").Value = "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.zh-cn.xlf"
$ws.Range("G4").Value = "2016-01-20 07:52:43"
$ws.Range("H4").Value = "Include"

# Row 5: .localization-config - unchanged
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# Rebuild hyperlinks in the new order, including the two brand-new ones
# (E4/F4) introduced by cb5280d3 now pointing at the completed handback.
$ws.Cells.Item(1, 1).Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5dc5f659cfefafe3787e5cc4efd5a0cb65c8d210/e2e/9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md", "", "", "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/91405156547e843225ebb942cd285fe03a3ed80e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.zh-cn.xlf", "", "", "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/91405156547e843225ebb942cd285fe03a3ed80e/e2e/9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md", "", "", "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/91405156547e843225ebb942cd285fe03a3ed80e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.zh-cn.xlf", "", "", "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.zh-cn.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/481cee7df47d66e29231fd5af014ae7c5826663a/e2e/ae38855b-7f4f-4934-8340-95b809897df5.md", "", "", "ae38855b-7f4f-4934-8340-95b809897df5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6db815d0008b55f2ec3ad2a556fc397339b02f63/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.zh-cn.xlf", "", "", "ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/eda1794def1d9d9bbf1bbe20c2e9e32e63863e20/e2e/ae38855b-7f4f-4934-8340-95b809897df5.md", "", "", "ae38855b-7f4f-4934-8340-95b809897df5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a9a59979f014c4d27a1072f8c525bf385b33ee52/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.zh-cn.xlf", "", "", "ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.zh-cn.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5dc5f659cfefafe3787e5cc4efd5a0cb65c8d210/e2e/cb5280d3-793c-4a18-8880-5345aa846a9e.md", "", "", "cb5280d3-793c-4a18-8880-5345aa846a9e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/91405156547e843225ebb942cd285fe03a3ed80e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.zh-cn.xlf", "", "", "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/91405156547e843225ebb942cd285fe03a3ed80e/e2e/9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md", "", "", "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/91405156547e843225ebb942cd285fe03a3ed80e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.zh-cn.xlf", "", "", "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.zh-cn.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/5dc5f659cfefafe3787e5cc4efd5a0cb65c8d210/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 2: 9d491b91 - now fully handed back
$ws.Range("A2").Value = "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.de-de.xlf"
$ws.Range("D2").Value = "2016-01-20 07:52:10"
$ws.Range("E2").Value = "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md"
$ws.Range("F2").Value = "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.de-de.xlf"
$ws.Range("G2").Value = "2016-01-20 07:53:02"
$ws.Range("H2").Value = "Include"

# Row 3: ae38855b - unchanged content, just shifted down a row
$ws.Range("A3").Value = "ae38855b-7f4f-4934-8340-95b809897df5.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.de-de.xlf"
$ws.Range("D3").Value = "2016-01-20 07:50:19"
$ws.Range("E3").Value = "ae38855b-7f4f-4934-8340-95b809897df5.md"
$ws.Range("F3").Value = "ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.de-de.xlf"
$ws.Range("G3").Value = "2016-01-20 07:51:13"
$ws.Range("H3").Value = "Include"

# Row 4: cb5280d3 dependency - now also handed back, mirroring 9d491b91's data
$ws.Range("A4").Value = "cb5280d3-793c-4a18-8880-5345aa846a9e.md"
$ws.Range("B4").Value = "Handed back: in sync with en-US"
$ws.Range("C4").Value = "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.de-de.xlf"
$ws.Range("D4").Value = "2016-01-20 07:52:10"
$ws.Range("E4").Value = "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md"
$ws.Range("F4").Value = "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.de-de.xlf"
$ws.Range("G4").Value = "2016-01-20 07:53:02"
$ws.Range("H4").Value = "Include"

# Row 5: .localization-config - unchanged
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# Rebuild hyperlinks in the new order, including the two brand-new ones
# (E4/F4) introduced by cb5280d3 now pointing at the completed handback.
$ws.Cells.Item(1, 1).Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5dc5f659cfefafe3787e5cc4efd5a0cb65c8d210/e2e/9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md", "", "", "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eb5d9fd32fe6273bf6ff6435a90ecddb22a77560/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.de-de.xlf", "", "", "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/eb5d9fd32fe6273bf6ff6435a90ecddb22a77560/e2e/9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md", "", "", "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/eb5d9fd32fe6273bf6ff6435a90ecddb22a77560/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.de-de.xlf", "", "", "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.de-de.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/481cee7df47d66e29231fd5af014ae7c5826663a/e2e/ae38855b-7f4f-4934-8340-95b809897df5.md", "", "", "ae38855b-7f4f-4934-8340-95b809897df5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2b761fb9cfa4d6e559bbbf4f63117e289686697d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.de-de.xlf", "", "", "ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/dfba68575d9ceb85fa982f0741a4894d1aee7ac7/e2e/ae38855b-7f4f-4934-8340-95b809897df5.md", "", "", "ae38855b-7f4f-4934-8340-95b809897df5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/059bd61d218349e60253cbe6f7d9b9b71d109bb9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.de-de.xlf", "", "", "ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.de-de.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/5dc5f659cfefafe3787e5cc4efd5a0cb65c8d210/e2e/cb5280d3-793c-4a18-8880-5345aa846a9e.md", "", "", "cb5280d3-793c-4a18-8880-5345aa846a9e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eb5d9fd32fe6273bf6ff6435a90ecddb22a77560/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.de-de.xlf", "", "", "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/eb5d9fd32fe6273bf6ff6435a90ecddb22a77560/e2e/9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md", "", "", "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/eb5d9fd32fe6273bf6ff6435a90ecddb22a77560/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.de-de.xlf", "", "", "9d491b91-1b9f-4b8b-afa4-49d9a35f8c4f.317ff1f062e3e538dcf689f199d31afb3d915e9f.de-de.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/5dc5f659cfefafe3787e5cc4efd5a0cb65c8d210/.localization-config", "", "", ".localization-config") | Out-Null
